# feat: add 2022-Q3 data
#
# The workbook has 3 sheets: "总计" (summary), "2022-Q2", "2022-Q1".
# We add a new "2022-Q3" sheet (a copy of the "2022-Q2" sheet, with
# refreshed figures) positioned right after "总计" and before "2022-Q2",
# and refresh the "总计" summary sheet so it lists all three quarters.

$wb = $excel.ActiveWorkbook

# --- 1. Create the "2022-Q3" sheet as a duplicate of "2022-Q2" ---------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)                       # new copy is inserted right before $q2
$q3 = $wb.Worksheets.Item(2)        # "总计" is #1, the new copy is now #2
$q3.Name = "2022-Q3"

# Refresh the Q3 numbers (fund-size, stock position, weight, market value
# and the rank all shifted with the new quarter's data).
$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "27.03"
$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "99.07"
$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "2.54"
$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "0.6866"
$q3.Range("H2").Value = 10

# --- 2. Refresh the "总计" summary sheet -------------------------------
$total = $wb.Worksheets.Item("总计")

# Add a 4th row (copy row 3's formatting so the new row matches the rest
# of the table), then shift every quarter label/value down by one row.
$total.Range("A3:D3").Copy()
$total.Range("A4:D4").PasteSpecial(-4122)   # xlPasteFormats

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.67

$total.Range("B3").Value = "2022-Q2"
$total.Range("D3").Value = 0.72

$total.Range("B2").Value = "2022-Q3"
$total.Range("D2").Value = 0.69

# Keep the same sheet active as before the edit.
$total.Activate()
